$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its current location (end of the
#    "Legnagyobb problema..." paragraph). It will be re-added later at the
#    end of the newly inserted content.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Move the existing "Metaprogramozasrol altalaban" .. "Rascal" block (7
#    paragraphs) so that it sits right after the "Legnagyobb problema..."
#    paragraph instead of after it. Using Cut+Paste (rather than retyping)
#    preserves the original paragraph formatting (pStyle/ilvl) and run
#    structure (e.g. the ".NET" / " Framework" two-run split) exactly.
# ---------------------------------------------------------------------------
$anchorIndex = 3   # "Legnagyobb problema..." paragraph; grows by 1 each iteration
for ($k = 0; $k -lt 7; $k++) {
    $srcIndex = $anchorIndex + 1
    $srcPara = $d.Paragraphs.Item($srcIndex)
    $srcPara.Range.Cut()

    $anchor = $d.Paragraphs.Item($anchorIndex)
    $destRange = $d.Range($anchor.Range.End, $anchor.Range.End)
    $destRange.Paste()

    $anchorIndex = $anchorIndex + 1
}
# $anchorIndex now points at the last paragraph of the moved block ("Rascal...").

# ---------------------------------------------------------------------------
# 3. Insert the brand new paragraphs about the advantages of metaprogramming
#    right after the moved block.
# ---------------------------------------------------------------------------
$newTexts = @(
    @{ Text = "Metaprogramozás előnyei/lehetőségei"; Level = 1 },
    @{ Text = "Fordító optimalizációja: inline direktíva, végrekurzió stb."; Level = 2 },
    @{ Text = "Deklaratív programozás támogatása"; Level = 2 },
    @{ Text = "Mintaillesztés implementációja (úgy, mint a Scala-ban) – nem hiszem, hogy ezt meg lehetne oldani"; Level = 2 },
    @{ Text = "Új paradigmák bevezetése: pl.: Design by Contract (elő- utófeltétel, invariánsok stb.)"; Level = 2 }
)

$prevParaIndex = $anchorIndex
foreach ($item in $newTexts) {
    $prevPara = $d.Paragraphs.Item($prevParaIndex)
    $prevPara.Range.InsertParagraphAfter()
    $newParaIndex = $prevParaIndex + 1
    $newPara = $d.Paragraphs.Item($newParaIndex)
    $newPara.Range.Text = $item.Text
    $newPara2 = $d.Paragraphs.Item($newParaIndex)
    $newPara2.Range.ListFormat.ListLevelNumber = $item.Level
    $prevParaIndex = $newParaIndex
}

# $prevParaIndex now points at the last inserted paragraph
# ("Uj paradigmak bevezetese..."). Re-add the _GoBack bookmark at the very
# end of its text, immediately before the paragraph mark.
#
# NOTE: placing a zero-length bookmark exactly at (paragraph.End - 1) is
# mishandled when that position is directly adjacent to the paragraph mark,
# so a temporary placeholder character is used to work around it: insert a
# throw-away character, anchor the bookmark just before it, then delete the
# placeholder again.
$lastNewPara = $d.Paragraphs.Item($prevParaIndex)
$tmpRange = $lastNewPara.Range
$tmpRange.Collapse(0)
$tmpRange.MoveEnd(1, -1) | Out-Null
$tmpRange.Collapse(0)
$tmpRange.InsertAfter("X")

$lastNewPara2 = $d.Paragraphs.Item($prevParaIndex)
$bookmarkPos = $lastNewPara2.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$lastNewPara3 = $d.Paragraphs.Item($prevParaIndex)
$placeholderPos = $lastNewPara3.Range.End - 2
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
